# Applies the "cryptos list" refresh described in the commit message:
#   - Price (column D) and Volume(1h) (column E) values are updated for every
#     coin row to the newly scraped figures.
#   - Three coins were re-ranked (rows 35-37: Kaspa, EthereumClassic, dogwifhat)
#     and two coins swapped places (rows 50-51: OKB, ImmutableX), so their Coin
#     name/Link/Price/Volume cells were replaced with the coin that now sits on
#     that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '90.347.87'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.83%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.096.88'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.25%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.25'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.97%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '621.88'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.21%  '

# Row 7
$ws.Range('E7').Value = '  -6.54%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.366'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.09%  '

# Row 9
$ws.Range('E9').Value = '  +0.17%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.099.39'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.04%  '

# Row 11
$ws.Range('E11').Value = '  -2.97%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.197'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.21%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.51%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.72'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.72%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.46'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.98%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.240.19'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.42%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.672.37'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.19%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.91'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +6.07%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.104.28'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.56%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000218'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.72%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.99'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.48%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '436.44'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.81%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.55'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.35%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.43%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.86'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.72%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.59'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.11%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.26'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.15%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.04'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.99%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.270.93'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.82%  '

# Row 30
$ws.Range('E30').Value = '  +0.00%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.28'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.33%  '

# Row 32
$ws.Range('E32').Value = '  -1.46%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.26%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.194'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.42%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.153'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +7.28%  '

# Row 36
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '25.79'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.56%  '

# Row 37
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.80'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.99%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.16'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.82%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '502.14'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.61%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.90'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.24%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.28'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.95%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0886'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.40%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.16'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.39%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.405'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.40%  '

# Row 45
$ws.Range('E45').Value = '  +0.00%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +50.91%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.90'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.81%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.687'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.04%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '151.02'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.27%  '

# Row 50
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.56'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.88%  '

# Row 51
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.33'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.43%  '
